# Automatic update of files.
# Increment the "Förändrad" (Changed/Updated) date in column C from
# 2025-04-17 (serial 45764) to 2025-04-18 (serial 45765) for all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45764) {
        $cell.Value2 = 45765
    }
}
